# fix: update LaporanPenjualanExport and related views to include location
# and format currency (point 19)
#
# The exported "Laporan Penjualan" template gains two blank rows above the
# table header so the header block (shop name / address) has room for the
# additional "location" line. Concretely: the header row that used to sit
# on row 5 (Tanggal / Nama Pelanggan / Kode Bbarang / ... / Total) moves
# down to row 7, leaving rows 4-6 blank beneath the existing title rows
# (1-2 merged "BDR HALL", row 3 "Jl. Tinumbu No.20 ...").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows directly above the header row (row 5), pushing
# the header row (and everything below it) down to row 7. Excel shifts all
# row-anchored content (values, styles, merged cells below the insert
# point) down automatically, which is exactly the A5->A7 ... L5->L7 move
# seen in the diff.
$ws.Rows("5:6").Insert() | Out-Null

# Leave the sheet selection on the (now empty) row 5, matching the
# saved view state in the edited workbook.
$ws.Rows("5:5").Select() | Out-Null
